# Auto-generated Excel COM-interop edit script
# Commit: Add data for 2025-09-30
# Applies targeted cell value updates across 21 worksheets in the
# 'cta-violent-crime-ytd' workbook, matching the authoritative diff.

$wb = $excel.ActiveWorkbook

# --- Citywide Totals ---
$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("C2").Value = 55
$ws.Range("G2").Value = 69
$ws.Range("J2").Value = 90
$ws.Range("E3").Value = 107
$ws.Range("F3").Value = 100
$ws.Range("I3").Value = 165
$ws.Range("J3").Value = 172
$ws.Range("B6").Value = 301
$ws.Range("C6").Value = 370
$ws.Range("D6").Value = 327
$ws.Range("E6").Value = 346
$ws.Range("F6").Value = 409
$ws.Range("G6").Value = 380
$ws.Range("H6").Value = 358
$ws.Range("I6").Value = 409
$ws.Range("J6").Value = 322
$ws.Range("L6").Value = 357
$ws.Range("B7").Value = 406
$ws.Range("C7").Value = 499
$ws.Range("D7").Value = 511
$ws.Range("E7").Value = 517
$ws.Range("F7").Value = 583
$ws.Range("G7").Value = 558
$ws.Range("H7").Value = 565
$ws.Range("I7").Value = 688
$ws.Range("J7").Value = 601
$ws.Range("L7").Value = 676

# --- Garfield Park ---
$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J3").Value = 11
$ws.Range("J7").Value = 35

# --- Grand Crossing ---
$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("B6").Value = 19
$ws.Range("B7").Value = 25

# --- Armour Square ---
$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("E3").Value = 4
$ws.Range("I5").Value = 11
$ws.Range("L5").Value = 15
$ws.Range("E6").Value = 8
$ws.Range("I6").Value = 15
$ws.Range("L6").Value = 21

# --- Humboldt Park ---
$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J5").Value = 9
$ws.Range("J6").Value = 15

# --- Englewood ---
$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("B6").Value = 27
$ws.Range("E6").Value = 18
$ws.Range("L6").Value = 33
$ws.Range("B7").Value = 30
$ws.Range("E7").Value = 27
$ws.Range("L7").Value = 55

# --- By Neighborhood ---
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("E5").Value = 8
$ws.Range("I5").Value = 15
$ws.Range("L5").Value = 21
$ws.Range("H6").Value = 2
$ws.Range("D8").Value = 26
$ws.Range("J8").Value = 36
$ws.Range("C18").Value = 2
$ws.Range("B27").Value = 2
$ws.Range("B28").Value = 30
$ws.Range("E28").Value = 27
$ws.Range("L28").Value = 55
$ws.Range("J32").Value = 35
$ws.Range("B36").Value = 25
$ws.Range("J41").Value = 15
$ws.Range("C47").Value = 21
$ws.Range("J47").Value = 10
$ws.Range("I53").Value = 109
$ws.Range("F56").Value = 2
$ws.Range("C65").Value = 16
$ws.Range("H77").Value = 22
$ws.Range("E78").Value = 6
$ws.Range("G87").Value = 5
$ws.Range("G94").Value = 4
$ws.Range("F97").Value = 6
$ws.Range("B98").Value = 406
$ws.Range("C98").Value = 499
$ws.Range("D98").Value = 511
$ws.Range("E98").Value = 517
$ws.Range("F98").Value = 583
$ws.Range("G98").Value = 558
$ws.Range("H98").Value = 565
$ws.Range("I98").Value = 688
$ws.Range("J98").Value = 601
$ws.Range("L98").Value = 676

# --- Rush & Division ---
$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("E3").Value = 2
$ws.Range("E5").Value = 6

# --- Loop ---
$ws = $wb.Worksheets.Item('Loop')
$ws.Range("I3").Value = 27
$ws.Range("I7").Value = 109

# --- Grand Boulevard ---
$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("C2").Value = 2

# --- North Lawndale ---
$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("C6").Value = 16

# --- Washington Heights ---
$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("G4").Value = 4
$ws.Range("G5").Value = 5

# --- Roseland ---
$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("H6").Value = 13
$ws.Range("H7").Value = 22

# --- Calumet Heights ---
$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("C4").Value = 2
$ws.Range("C5").Value = 2

# --- Edgewater ---
$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("B4").Value = 2
$ws.Range("B5").Value = 2

# --- Wrigleyville ---
$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Range("F3").Value = 1
$ws.Range("F6").Value = 6

# --- Lake View ---
$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("J2").Value = 2
$ws.Range("C5").Value = 12
$ws.Range("C6").Value = 21
$ws.Range("J6").Value = 10

# --- Mckinley Park ---
$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("E5").Value = 1
$ws.Range("E6").Value = 2

# --- West Town ---
$ws = $wb.Worksheets.Item('West Town')
$ws.Range("G2").Value = 1
$ws.Range("G6").Value = 4

# --- Austin ---
$ws = $wb.Worksheets.Item('Austin')
$ws.Range("D5").Value = 18
$ws.Range("J5").Value = 18
$ws.Range("D6").Value = 26
$ws.Range("J6").Value = 36

# --- Ashburn ---
$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("F4").Value = 1
$ws.Range("F5").Value = 2

